$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently shows the shared string "R40". The edit replaces its
# displayed value with the text "1" -- a *text* value (not a number), so
# it must land back in the workbook as a shared string, not as numeric 1.
#
# A plain `Range.Value = "1"` gets auto-coerced to the number 1 by Excel's
# usual "looks like a number" heuristic, and forcing text via a leading
# apostrophe or via NumberFormat="@" both stamp a new/changed style onto
# the cell. Neither preserves cell B11's existing style untouched.
#
# Instead: build the literal text "1" as a formula result in a scratch
# cell, copy it, and paste *values only* into B11. PasteSpecial(values)
# carries over the source's text data type without touching B11's
# existing number format/style, and leaves no trace behind once the
# scratch cell is cleared.
$scratch = $ws.Range("ZZ1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false
